$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 121
$ws.Range("B6").Value = 0.46552083333333333
$ws.Range("C6").Value = 0.25

$ws.Range("A7").Value = 121
$ws.Range("B7").Value = 0.46249999999999997
$ws.Range("C7").Value = 0.4

$ws.Range("B6:B7").NumberFormat = $ws.Range("B5").NumberFormat

$ws.Range("A8").Select()
